# Applies the "OpcionesSist" model change: adds a new "Codigo" field row,
# restyles / merges the title row, colors the sheet tab, fixes selections,
# and moves the active sheet/tab selection to "4. OpcionesSist".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Sheet "3. Usuarios": it used to be the active/selected tab. We leave
#    its own selection (K5) untouched; tabSelected will move away from it
#    automatically once another sheet becomes active later in this script.
# ---------------------------------------------------------------------
$wsUsuarios = $wb.Worksheets.Item("3. Usuarios")

# ---------------------------------------------------------------------
# 2) Sheet "5. PermisosRoles": selection changes from A2:XFD2 to G10.
#    Do this BEFORE activating "4. OpcionesSist" so the final active tab
#    ends up being OpcionesSist, not this one.
# ---------------------------------------------------------------------
$wsPermisos = $wb.Worksheets.Item("5. PermisosRoles")
$wsPermisos.Range("G10").Select()

# ---------------------------------------------------------------------
# 3) Sheet "4. OpcionesSist": main set of changes.
# ---------------------------------------------------------------------
$wsOpciones = $wb.Worksheets.Item("4. OpcionesSist")

# Tab color (red)
$wsOpciones.Tab.Color = 192

# Insert a new row 4 ("Codigo" / "String" / 30) pushing the former rows
# 4 and 5 down to 5 and 6.
$wsOpciones.Rows("4:4").Insert()
$wsOpciones.Range("A4").Value = "Codigo"
$wsOpciones.Range("B4").Value = "String"
$wsOpciones.Range("C4").Value = 30

# Column A width (best-fit style width for the new content)
$wsOpciones.Columns("A:A").ColumnWidth = 10.5

# Title row formatting: bold, size 14 font, left aligned, merged across A1:C1
$wsOpciones.Range("A1:C1").Font.Bold = $true
$wsOpciones.Range("A1:C1").Font.Size = 14
$wsOpciones.Range("A1:C1").HorizontalAlignment = -4131
$wsOpciones.Range("A1:C1").Merge()

# Select A5 and make this the active sheet/tab (this also sets
# tabSelected="1" on this sheet and clears it from whichever sheet had it
# before, and updates workbook activeTab).
$wsOpciones.Range("A5").Select()

Write-Host "Edit complete"
